$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns B/C/E hold values that Excel will never misread as numbers
# (plain text, or percentages padded with spaces), so a direct .Value
# assignment is sufficient and keeps the cell unstyled, as in the source.
#
# Column D holds price figures. Some look like plain decimals (e.g. "11.30")
# which Excel.Value would silently coerce to the Number 11.3, losing the
# trailing zero and the original text type. To keep these as literal text
# (matching the workbook's inlineStr cells) we force NumberFormat="@" before
# the assignment, then ClearFormats() right after so no stray cell style
# (other than the pre-existing "no style") is left behind.

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '69.302.83'
$ws.Range('D2').ClearFormats()
$ws.Range('E2').Value = '  -4.34%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.500.18'
$ws.Range('D3').ClearFormats()
$ws.Range('E3').Value = '  -5.55%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '575.06'
$ws.Range('D5').ClearFormats()
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '165.79'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -5.00%  '
$ws.Range('E7').Value = '  +0.06%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.513'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  -1.03%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '2.498.99'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  -5.55%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.155'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -9.80%  '
$ws.Range('E12').Value = '  -4.08%  '
$ws.Range('E13').Value = '  -2.36%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '2.960.31'
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -5.43%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '69.332.31'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -4.07%  '
$ws.Range('E16').Value = '  -7.38%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '24.63'
$ws.Range('D17').ClearFormats()
$ws.Range('E17').Value = '  -4.99%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '2.495.49'
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -6.02%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.30'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  -6.96%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '7.72'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.09%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '345.79'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -7.16%  '
$ws.Range('E22').Value = '  -5.87%  '
$ws.Range('E23').Value = '  +0.11%  '
$ws.Range('E24').Value = '  -6.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.11'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -4.01%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.93'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -7.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '8.84'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -8.15%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.629.95'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -5.45%  '
$ws.Range('E29').Value = '  -0.68%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.0₃0890'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  -7.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.81'
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -2.69%  '
$ws.Range('B32').Value = 'Fetch.AI'
$ws.Range('C32').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.25'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -2.54%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '454.58'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -8.64%  '
$ws.Range('E34').Value = '  -3.60%  '
$ws.Range('E35').Value = '  -0.02%  '
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '153.24'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -5.12%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '18.95'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +0.33%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.30'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -5.13%  '
$ws.Range('E40').Value = '  +0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.70'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -4.37%  '
$ws.Range('E43').Value = '  -8.97%  '
$ws.Range('E44').Value = '  -15.32%  '
$ws.Range('E45').Value = '  -2.67%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.27'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -12.02%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '142.22'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -6.93%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.523'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -4.68%  '
$ws.Range('E49').Value = '  -5.05%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.58'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -5.70%  '
$ws.Range('E51').Value = '  -2.58%  '
